$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '62.827.86'
$ws.Range("E2").Value = '  -0.62%  '
$ws.Range("D3").Value = '2.465.41'
$ws.Range("E3").Value = '  -0.53%  '
$ws.Range("E4").Value = '  -0.02%  '
$ws.Range("D5").Value = "'572.08"
$ws.Range("E5").Value = '  -0.78%  '
$ws.Range("D6").Value = "'147.44"
$ws.Range("E7").Value = '  +0.00%  '
$ws.Range("D8").Value = "'0.532"
$ws.Range("E8").Value = '  -1.42%  '
$ws.Range("E9").Value = '  +0.03%  '
$ws.Range("D10").Value = "'0.163"
$ws.Range("E10").Value = '  -0.04%  '
$ws.Range("E11").Value = '  -1.21%  '
$ws.Range("E12").Value = '  -1.51%  '
$ws.Range("D13").Value = "'29.19"
$ws.Range("E13").Value = '  +2.16%  '
$ws.Range("D15").Value = '2.908.77'
$ws.Range("E15").Value = '  -0.68%  '
$ws.Range("D16").Value = '62.735.00'
$ws.Range("E16").Value = '  -0.66%  '
$ws.Range("D17").Value = '2.464.96'
$ws.Range("E17").Value = '  -0.56%  '
$ws.Range("E18").Value = '  -5.83%  '
$ws.Range("D19").Value = "'10.73"
$ws.Range("E19").Value = '  -2.55%  '
$ws.Range("E20").Value = '  +4.65%  '
$ws.Range("D21").Value = "'4.16"
$ws.Range("E21").Value = '  +0.50%  '
$ws.Range("D22").Value = "'321.45"
$ws.Range("E22").Value = '  -2.59%  '
$ws.Range("D24").Value = "'10.24"
$ws.Range("E24").Value = '  +3.82%  '
$ws.Range("D25").Value = "'64.79"
$ws.Range("E25").Value = '  -2.13%  '
$ws.Range("D26").Value = "'643.07"
$ws.Range("E26").Value = '  -1.99%  '
$ws.Range("E27").Value = '  -0.63%  '
$ws.Range("D28").Value = '0.0₃0964'
$ws.Range("E28").Value = '  -2.69%  '
$ws.Range("E29").Value = '  +0.82%  '
$ws.Range("E30").Value = '  -3.95%  '
$ws.Range("E31").Value = '  -2.02%  '
$ws.Range("E32").Value = '  -2.06%  '
$ws.Range("E33").Value = '  -0.01%  '
$ws.Range("E34").Value = '  -0.09%  '
$ws.Range("E35").Value = '  -3.13%  '
$ws.Range("E36").Value = '  -1.95%  '
$ws.Range("E37").Value = '  -1.36%  '
$ws.Range("E38").Value = '  -1.67%  '
$ws.Range("D39").Value = "'18.55"
$ws.Range("E39").Value = '  -1.17%  '
$ws.Range("D40").Value = "'149.80"
$ws.Range("E40").Value = '  -0.30%  '
$ws.Range("D41").Value = "'2.64"
$ws.Range("E41").Value = '  -1.36%  '
$ws.Range("D42").Value = "'1.73"
$ws.Range("E42").Value = '  -1.66%  '
$ws.Range("B43").Value = 'BabyDogeCoin'
$ws.Range("C43").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D43").Value = '0.0₆0307'
$ws.Range("E43").Value = '  -4.08%  '
$ws.Range("B44").Value = 'USDe'
$ws.Range("C44").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("D44").Value = "'0.999"
$ws.Range("E44").Value = '  +0.03%  '
$ws.Range("D45").Value = "'154.08"
$ws.Range("E45").Value = '  -0.63%  '
$ws.Range("D46").Value = "'15.39"
$ws.Range("E46").Value = '  +0.97%  '
$ws.Range("D48").Value = "'20.27"
$ws.Range("E48").Value = '  -0.74%  '
$ws.Range("E49").Value = '  -0.39%  '
$ws.Range("D50").Value = "'0.0510"
$ws.Range("E50").Value = '  -1.12%  '
$ws.Range("E51").Value = '  -1.66%  '
